# Updates the cryptocurrency price (column D) and volume-change (column E) values
# on the single worksheet, matching the latest scrape from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($range, $text) {
    # Force the cell to hold the given string as text, even when the text
    # looks like a number (e.g. "528.20"), by using a leading apostrophe to
    # suppress Excels automatic number conversion, then clearing the
    # resulting quote-prefix formatting so the cell style is left untouched.
    $range.Value = "'" + $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "60.639.65"
$ws.Range("E2").Value = "  -1.66%  "
Set-TextValue $ws.Range("D3") "2.902.98"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "528.20"
$ws.Range("E5").Value = "  -2.44%  "
Set-TextValue $ws.Range("D6") "144.23"
$ws.Range("E6").Value = "  -5.64%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.31%  "
Set-TextValue $ws.Range("D9") "2.910.76"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("E10").Value = "  -3.50%  "
Set-TextValue $ws.Range("D11") "6.03"
$ws.Range("E11").Value = "  -1.16%  "
Set-TextValue $ws.Range("D12") "0.363"
$ws.Range("E12").Value = "  -0.53%  "
Set-TextValue $ws.Range("D13") "3.408.38"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("E14").Value = "  +2.59%  "
Set-TextValue $ws.Range("D15") "60.630.52"
$ws.Range("E15").Value = "  -1.75%  "
Set-TextValue $ws.Range("D16") "22.79"
$ws.Range("E16").Value = "  -3.68%  "
Set-TextValue $ws.Range("D17") "2.906.21"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("E21").Value = "  -4.58%  "
Set-TextValue $ws.Range("D22") "6.64"
$ws.Range("E22").Value = "  -0.16%  "
Set-TextValue $ws.Range("D23") "0.999"
$ws.Range("E23").Value = "  -0.06%  "
Set-TextValue $ws.Range("D24") "5.68"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("E27").Value = "  -2.85%  "
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  -4.27%  "
Set-TextValue $ws.Range("D30") "0.0₃0861"
$ws.Range("E30").Value = "  -7.48%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("E33").Value = "  -3.55%  "
Set-TextValue $ws.Range("D34") "151.50"
$ws.Range("E34").Value = "  -4.60%  "
Set-TextValue $ws.Range("D35") "4.39"
$ws.Range("E35").Value = "  -5.32%  "
Set-TextValue $ws.Range("D36") "5.57"
$ws.Range("E36").Value = "  -6.53%  "
$ws.Range("E37").Value = "  -4.58%  "
$ws.Range("E38").Value = "  -5.47%  "
Set-TextValue $ws.Range("D39") "37.67"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  -4.36%  "
Set-TextValue $ws.Range("D41") "3.73"
$ws.Range("E41").Value = "  -4.81%  "
Set-TextValue $ws.Range("D42") "2.291.28"
$ws.Range("E42").Value = "  -4.95%  "
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("E45").Value = "  -7.13%  "
Set-TextValue $ws.Range("D46") "0.997"
$ws.Range("E46").Value = "  +0.00%  "
Set-TextValue $ws.Range("D47") "5.02"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("E50").Value = "  -2.01%  "
Set-TextValue $ws.Range("D51") "251.89"
$ws.Range("E51").Value = "  -5.16%  "
